$wb = $excel.ActiveWorkbook

# --- Duplicate the "Non-Stopping TTO" sheet into a new "New TTO" sheet ---
$src = $wb.Worksheets.Item("Non-Stopping TTO")
$src.Copy([System.Reflection.Missing]::Value, $src)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "New TTO"

# The copied "F" column (questionid) values need to be bumped from 5 to 6
# on the new sheet, representing the next questionnaire version.
$newSheet.Range("F2:F50").Value = 6

# The previously active sheet ("TTO") had its selection updated too.
$wsTto = $wb.Worksheets.Item("TTO")
$wsTto.Activate()
$wsTto.Range("F2").Select()

# View state tweaks that came along with the edit:
# the new sheet becomes the active / selected tab, scrolled down a bit,
# keeping the same selected range as the source sheet.
$newSheet.Activate()
$newSheet.Range("F2:F50").Select()
